# bets-2023-2.xlsx — "Add files via upload"
# Restructures the "resumen" sheet from a 2-column (label/value) layout
# into a wide table with a TIPO/ counter column, a data row and a row of
# array formulas that apply the latest percentage from bets!$M:$M, and
# flips which sheet/tab is active.

$wb = $excel.ActiveWorkbook

$betsWs    = $wb.Worksheets.Item("bets")
$resumenWs = $wb.Worksheets.Item("resumen")

# ---------------------------------------------------------------------
# 1) "resumen" sheet: rebuild the old A1:B5 label/value table into the
#    new A1:F3 wide table.
# ---------------------------------------------------------------------

# Wipe the old contents (keeps column B's number-format style) then drop
# the now-unused rows 4:5 so the sheet's used range shrinks back down.
$resumenWs.Range("A1:B5").ClearContents()
$resumenWs.Rows("4:5").Delete()

# Row 1: header labels
$resumenWs.Range("A1").Value = "TIPO"
$resumenWs.Range("B1").Value = "M"
$resumenWs.Range("C1").Value = "V"
$resumenWs.Range("D1").Value = "E"
$resumenWs.Range("E1").Value = "C"
$resumenWs.Range("F1").Value = "CC"

# Row 2: the original figures, now spread across the row instead of down
# a column.
$resumenWs.Range("A2").Value = 1
$resumenWs.Range("B2").Value = 462440
$resumenWs.Range("C2").Value = 54753.06
$resumenWs.Range("D2").Value = 81628.56
$resumenWs.Range("E2").Value = 12896.82
$resumenWs.Range("F2").Value = 25793.64

# Give B2:F2 the "Millares" thousands-separator style, matching B1's
# pre-existing formatting.
$resumenWs.Range("B2:F2").NumberFormat = '_-* #,##0_-;\-* #,##0_-;_-* "-"_-;_-@_-'

# Row 3: counter + array formulas that grow each figure by the latest
# percentage recorded in bets!M (the last non-empty cell in that column).
$resumenWs.Range("A3").Value = 2

$resumenWs.Range("B3").FormulaArray = "=B2+B2*INDEX(bets!`$M:`$M,COUNTA(bets!`$M:`$M))/100"
$resumenWs.Range("C3").FormulaArray = "=C2+C2*INDEX(bets!`$M:`$M,COUNTA(bets!`$M:`$M))/100"
$resumenWs.Range("D3").FormulaArray = "=D2+D2*INDEX(bets!`$M:`$M,COUNTA(bets!`$M:`$M))/100"
$resumenWs.Range("E3").FormulaArray = "=E2+E2*INDEX(bets!`$M:`$M,COUNTA(bets!`$M:`$M))/100"
$resumenWs.Range("F3").FormulaArray = "=F2+F2*INDEX(bets!`$M:`$M,COUNTA(bets!`$M:`$M))/100"

$resumenWs.Range("B3:F3").NumberFormat = '_-* #,##0_-;\-* #,##0_-;_-* "-"_-;_-@_-'

# Selection inside "resumen" moves to G7.
$resumenWs.Range("G7").Select()

# ---------------------------------------------------------------------
# 2) "bets" sheet: selection moves from K22 to M15 and it loses the
#    "active tab" state (resumen becomes the active tab below).
# ---------------------------------------------------------------------
$betsWs.Range("M15").Select()

# ---------------------------------------------------------------------
# 3) Flip the active sheet/tab to "resumen" (workbook re-opens on it).
# ---------------------------------------------------------------------
$resumenWs.Activate()
